# Apply Betfair Back/Lay odds updates for 2026-02-11
# (values taken from the commit "Atualizando o arquivo XLSX")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("F6").Value = 3.7
$ws.Range("J6").Value = 3.1

# Row 7
$ws.Range("H7").Value = 2.26
$ws.Range("I7").Value = 2.42
$ws.Range("J7").Value = 3.45
$ws.Range("K7").Value = 3.7
$ws.Range("Q7").Value = 1.83

# Row 8
$ws.Range("P8").Value = 1.6

# Row 9
$ws.Range("J9").Value = 2.96
$ws.Range("P9").Value = 1.58

# Row 10
$ws.Range("F10").Value = 2.46
$ws.Range("H10").Value = 2.68
$ws.Range("I10").Value = 2.86
$ws.Range("K10").Value = 4.1
$ws.Range("P10").Value = 2.54
$ws.Range("Q10").Value = 1.48

# Row 11
$ws.Range("F11").Value = 1.69
$ws.Range("G11").Value = 1.75
$ws.Range("H11").Value = 4.4
$ws.Range("I11").Value = 6.2
$ws.Range("P11").Value = 2.6
$ws.Range("Q11").Value = 1.44

# Row 15
$ws.Range("P15").Value = 2.16
$ws.Range("Q15").Value = 1.74

# Row 17
$ws.Range("AB17").Value = 7.6
$ws.Range("AC17").Value = 9.6
$ws.Range("AD17").Value = 27
$ws.Range("AE17").Value = 1000
$ws.Range("AI17").Value = 120
$ws.Range("AJ17").Value = 14.5
$ws.Range("AM17").Value = 160
$ws.Range("H17").Value = 6.8
$ws.Range("I17").Value = 7.2
$ws.Range("J17").Value = 4.2
$ws.Range("N17").Value = 3.75
$ws.Range("P17").Value = 1.95
$ws.Range("Q17").Value = 2.02
$ws.Range("R17").Value = 1.35
$ws.Range("S17").Value = 3.65
$ws.Range("U17").Value = 1.89
$ws.Range("X17").Value = 14.5
$ws.Range("Z17").Value = 60

# Row 18
$ws.Range("AA18").Value = 80
$ws.Range("AD18").Value = 16.5
$ws.Range("AE18").Value = 46
$ws.Range("AJ18").Value = 24
$ws.Range("AM18").Value = 80
$ws.Range("AO18").Value = 42
$ws.Range("F18").Value = 2
$ws.Range("G18").Value = 2.02
$ws.Range("Q18").Value = 1.82
$ws.Range("X18").Value = 17
$ws.Range("Y18").Value = 17

# Row 19
$ws.Range("AA19").Value = 180
$ws.Range("AB19").Value = 8.4
$ws.Range("AD19").Value = 21
$ws.Range("AE19").Value = 80
$ws.Range("AI19").Value = 80
$ws.Range("AL19").Value = 38
$ws.Range("AM19").Value = 120
$ws.Range("AN19").Value = 11.5
$ws.Range("F19").Value = 1.78
$ws.Range("I19").Value = 5.6
$ws.Range("M19").Value = 1.08
$ws.Range("N19").Value = 3.75
$ws.Range("O19").Value = 1.34
$ws.Range("P19").Value = 1.94
$ws.Range("Q19").Value = 2.02
$ws.Range("R19").Value = 1.36
$ws.Range("S19").Value = 3.65
$ws.Range("T19").Value = 1.96
$ws.Range("X19").Value = 14

# Row 20
$ws.Range("AF20").Value = 9.4
$ws.Range("AG20").Value = 10.5
$ws.Range("AH20").Value = 24
$ws.Range("AL20").Value = 29
$ws.Range("H20").Value = 9.199999999999999
$ws.Range("I20").Value = 9.6
$ws.Range("T20").Value = 1.87
$ws.Range("U20").Value = 2.12

# Row 22
$ws.Range("H22").Value = 5.3
$ws.Range("I22").Value = 6
$ws.Range("J22").Value = 3.8
$ws.Range("P22").Value = 1.89
$ws.Range("Q22").Value = 1.92

# Row 23
$ws.Range("F23").Value = 2.86
$ws.Range("H23").Value = 2.6
$ws.Range("J23").Value = 3.25

# Row 24
$ws.Range("J24").Value = 3.6
$ws.Range("P24").Value = 2.06

# Row 25
$ws.Range("AG25").Value = 14.5
$ws.Range("AL25").Value = 980
$ws.Range("H25").Value = 19.5
$ws.Range("K25").Value = 10
$ws.Range("T25").Value = 2.16
$ws.Range("U25").Value = 1.74
$ws.Range("Y25").Value = 75
$ws.Range("Z25").Value = 230

# Row 26
$ws.Range("F26").Value = 1.52
$ws.Range("G26").Value = 1.62
$ws.Range("H26").Value = 6.2
$ws.Range("I26").Value = 7.6
$ws.Range("J26").Value = 4.4
$ws.Range("P26").Value = 2.2
$ws.Range("Q26").Value = 1.68

# Row 28
$ws.Range("P28").Value = 2.38

# Row 29
$ws.Range("F29").Value = 2.66
$ws.Range("G29").Value = 3.25
$ws.Range("H29").Value = 2.72
$ws.Range("I29").Value = 3.4

# Row 30
$ws.Range("F30").Value = 1.95
$ws.Range("G30").Value = 2.16
$ws.Range("H30").Value = 3.85
$ws.Range("I30").Value = 5.1
$ws.Range("J30").Value = 3.45
$ws.Range("P30").Value = 1.87
$ws.Range("Q30").Value = 1.9

# Row 31
$ws.Range("G31").Value = 1.86
$ws.Range("H31").Value = 4.5
$ws.Range("J31").Value = 3.9
$ws.Range("K31").Value = 4.7
$ws.Range("P31").Value = 2.2
$ws.Range("Q31").Value = 1.65

# Row 32
$ws.Range("J32").Value = 3.85
$ws.Range("Q32").Value = 1.77

# Row 33
$ws.Range("AB33").Value = 19
$ws.Range("AG33").Value = 20
$ws.Range("AJ33").Value = 180
$ws.Range("AN33").Value = 85
$ws.Range("AO33").Value = 10.5
$ws.Range("I33").Value = 1.78
$ws.Range("N33").Value = 0
$ws.Range("O33").Value = 0
$ws.Range("Q33").Value = 1.88
$ws.Range("S33").Value = 3.25
$ws.Range("T33").Value = 1.85
$ws.Range("U33").Value = 2.12

# Row 34
$ws.Range("G34").Value = 2.54
$ws.Range("H34").Value = 3.25
$ws.Range("I34").Value = 3.65
$ws.Range("P34").Value = 1.74

# Row 35
$ws.Range("G35").Value = 2.64
$ws.Range("H35").Value = 3.4
$ws.Range("K35").Value = 3.2
$ws.Range("P35").Value = 1.54
$ws.Range("Q35").Value = 2.62

# Row 36
$ws.Range("P36").Value = 1.79
$ws.Range("Q36").Value = 2.1

# Row 37
$ws.Range("Q37").Value = 2.06
